$d = $word.ActiveDocument

$p2 = $d.Paragraphs.Item(12)
$r2 = $p2.Range
$r2.Find.Execute("Lin. (2020). Cross-Domain", $true, $false, $false, $false, $false, $true, 1, $false, "Lin. Cross-Domain", 2) | Out-Null

$p2b = $d.Paragraphs.Item(12)
$r2b = $p2b.Range
$appendText = "In IEEE TRANSACTIONS ON PATTERN ANALYSIS AND MACHINE INTELLIGENCE (TPAMI'21)."
$r2b.Find.Execute("Adversarial Graph Learning. ", $true, $false, $false, $false, $false, $true, 1, $false, ("Adversarial Graph Learning. " + $appendText), 2) | Out-Null

$p2c = $d.Paragraphs.Item(12)
Write-Output ("Para 12: " + $p2c.Range.Text)

$full = $p2c.Range
$appendStart = $full.End - 1 - $appendText.Length
$part1Len = "In IEEE TRANSACTIONS ON PATTERN ANALYSIS AND MACHINE INTELLIGENCE (TPAMI".Length
$part2Len = "'".Length
$rA = $d.Range($appendStart, $appendStart + $part1Len)
$rB = $d.Range($appendStart + $part1Len, $appendStart + $part1Len + $part2Len)
$rC = $d.Range($appendStart + $part1Len + $part2Len, $full.End - 1)
Write-Output ("rA: [" + $rA.Text + "]")
Write-Output ("rB: [" + $rB.Text + "]")
Write-Output ("rC: [" + $rC.Text + "]")
$rA.Font.Bold = 0
$rB.Font.Bold = 0
$rC.Font.Bold = 0
Write-Output "done"
